$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Add a new row (row 10) with:
#   A10 = 9
#   C10 = "Sind Sie über 60?"
#   D10 = "Ja:boolean"
# Only A10, C10 and D10 get values (B10 and E10 stay empty), matching
# the style used by the existing data rows (A9, C9, D9).

$ws.Range("A9").Copy() | Out-Null
$ws.Range("A10").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("C9").Copy() | Out-Null
$ws.Range("C10").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("D9").Copy() | Out-Null
$ws.Range("D10").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A10").Value = 9
# Set D before C so the new shared-string entries are appended in the
# same order as the target workbook (D10 -> "Ja:boolean" first).
$ws.Range("D10").Value = "Ja:boolean"
$ws.Range("C10").Value = "Sind Sie über 60?"

$wb.Save()
